$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime for file 475b1cd4... (row 2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-11 22:32:32"
$wsZhCn.Range("H2").Value = "2016-03-11 22:32:48"

# de-de sheet: update Correspond Handoff/Handback Datetime for file 475b1cd4... (row 2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-11 22:32:35"
$wsDeDe.Range("H2").Value = "2016-03-11 22:32:54"
